$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @("MCH342-1", "LEGAL- GUGULETHU CASES (TO BE SORTED)", "Series", "1 Box", "LOCATION: 33I | GRAP COUNT NUMER: NONE"),
    @("MCH342-2", "CHRIS KNIGHT TRIAL, AFFIDAVITS, C.A.J ODENDAAL CASES", "Series", "1 Box", "LOCATION: 33I | GRAP COUNT NUMER: NONE"),
    @("MCH342-3", "COURT CASES, COMPLAINT AGAINST SOUTH AFRICAN POLICE RE: ASHLEY KRIEL, GUGULETHU CASES, ATHLONE, MANENBERG, STATEMENTS BY GRANT FAHRENFORT, SHANTEK FICK, COLEEN FICK, YAZEED BAKER, MASSOR MOLLAGE, ESMAT NORDIEN, ESHAAM NORDIEN, BRIAN ARENDSE, CLEMENT JOHN MEYER, SWARTMAN JULIUS PAKAMA", "Series", "1 Box", "LOCATION: 33I | GRAP COUNT NUMER: NONE"),
    @("MCH342-4", "TROJAN HORSE CRIMINAL TRIAL (TO BE SORTED)", "Series", "1 Box", "LOCATION: 33I | GRAP COUNT NUMER: NONE"),
    @("MCH342-5", "VIDEOS- VIDEO EVIDENCE TO SOME CASES (TO BE SORTED)", "Series", "1 Box", "LOCATION: 33I | GRAP COUNT NUMER: NONE")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $r = $r + 1
}
